$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Professional summary paragraph: neutralize language
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed",
    2) | Out-Null

# -----------------------------------------------------------------
# 2) Siege Analytics bullet: split out "50M" into its own bold run
#    (matching the formatting of the adjacent "23%"/"64%" runs).
#    Scope the bold-ing to only the text just replaced, so the other
#    "50M" occurrences (summary / impact line) are left as plain text.
# -----------------------------------------------------------------
$bulletRng = $d.Content
$bulletRng.Find.Execute(
    "race coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "race coding errors affecting 50M voters, developed",
    2) | Out-Null

$bulletText = $bulletRng.Text
$idx50M = $bulletText.IndexOf("50M")
if ($idx50M -ge 0) {
    $subStart = $bulletRng.Start + $idx50M
    $subEnd = $subStart + 3
    $subRng = $d.Range($subStart, $subEnd)
    $subRng.Bold = 1
    $subRng.Font.Color = 5258796
}

# -----------------------------------------------------------------
# 3) Move "Field Director - The Feldman Group" block from its
#    current location (after "Programmer - Lake Research Partners",
#    right before "KEY PROJECTS") to before
#    "Software Engineer - Salsa Labs" (i.e. right after the last
#    "Research Director - PCCC" bullet).
# -----------------------------------------------------------------

function Get-ParaIndexByText($doc, $text) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.TrimEnd() -eq $text) {
            return $i
        }
    }
    return -1
}

# Step A: move just the Heading3 line alone (preserves its paragraph style)
$headingIdx = Get-ParaIndexByText $d "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012"
$pHeading = $d.Paragraphs.Item($headingIdx)
$rngHeading = $d.Range($pHeading.Range.Start, $pHeading.Range.End)
$rngHeading.Cut() | Out-Null

$lastPccBulletIdx = Get-ParaIndexByText $d "• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"
$pDest = $d.Paragraphs.Item($lastPccBulletIdx)
$insertPoint = $d.Range($pDest.Range.End, $pDest.Range.End)
$insertPoint.Paste() | Out-Null

# Step B: move the 3 body paragraphs + the "Political Campaign Management"
# sub-heading right after the heading we just relocated.
$bodyStartIdx = Get-ParaIndexByText $d "Political Campaign Management"
$bodyEndIdx = Get-ParaIndexByText $d "• Created custom reports and data visualizations based on specific client requirements"
$pBodyStart = $d.Paragraphs.Item($bodyStartIdx)
$pBodyEnd = $d.Paragraphs.Item($bodyEndIdx)
$rngBody = $d.Range($pBodyStart.Range.Start, $pBodyEnd.Range.End)
$rngBody.Cut() | Out-Null

$newHeadingIdx = Get-ParaIndexByText $d "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012"
$pNewHeading = $d.Paragraphs.Item($newHeadingIdx)
$insertPoint2 = $d.Range($pNewHeading.Range.End, $pNewHeading.Range.End)
$insertPoint2.Paste() | Out-Null

# -----------------------------------------------------------------
# 4) Key Projects impact line: neutralize language
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved electoral prediction accuracy by 22%",
    2) | Out-Null

Write-Output "done"
